$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new row at row 8 (pushes existing rows 8+ down by one,
# formulas referencing shifted cells are adjusted automatically).
$ws.Rows.Item(8).Insert()

# Carry the formatting of the row above (row 7) into the newly
# inserted blank row 8, matching the look of the surrounding entries.
$ws.Range("A7:F7").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row with the uploaded ledger entry.
$ws.Range("B8").Value = 45358
$ws.Range("C8").Value = "b23-24MQ408"
$ws.Range("D8").Value = "Putzmeister Concrete Machines Pvt Ltd"
$ws.Range("E8").Value = 168741

# Row 7 no longer carries the running-total formula; it now lives on
# the newly-inserted row 8 and includes the new entry.
$ws.Range("F7").ClearContents()
$ws.Range("F8").Formula = "=E5+E6+E7+E8"

# Sheet2 ("Sale 22-23") becomes the active tab/sheet, with F7 selected.
$ws.Activate()
$ws.Range("F7").Select()
